# feat: add 2022-Q1 data
#
# The workbook contains one fund-holding sheet per quarter (2020-Q4 ..
# 2021-Q4) plus a running "总计" (grand total) summary sheet at the end.
# This change:
#   1. Adds a new "2022-Q1" fund-holding sheet (in the same position the
#      summary sheet used to occupy).
#   2. Moves the "总计" summary sheet one slot later (after "2022-Q1") and
#      adds a new leading row to it for the 2022-Q1 totals.
#
# To end up with the same sheetId / r:id numbering that the original
# authoring tool produced, we reuse the existing (last) worksheet - which
# is currently named "总计" - as the new "2022-Q1" sheet, and create the
# new "总计" sheet by duplicating it (so it starts out with the same
# look & feel / cell styles) immediately afterwards, then replace its
# contents with the refreshed summary table.

$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$q1 = $wb.Worksheets.Item($lastIndex)          # currently named "总计"

# Duplicate it right after itself -> this clone becomes the new "总计"
# sheet. Doing it now (before we touch $q1) means both sheets start out
# with identical cell styles (s="2" header/index cells etc.) that we can
# keep reusing locally via copy/paste-format, without ever having to copy
# formatting across different sheets.
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)    # the clone

$q1.Name = "2022-Q1"
$total.Name = "总计"

# ---------------------------------------------------------------------------
# Step 1: rewrite $q1 ("2022-Q1") as a fund-holding table
# ---------------------------------------------------------------------------

# The sheet currently has a 4-column (A:D) header/data layout; extend the
# styled header formatting (s="2") from D1 over the three new header cells.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Extend the styled index column (s="2") from A6 down to A12, for the
# extra fund rows.
$q1.Range("A6").Copy()
$q1.Range("A7:A12").PasteSpecial(-4122)

$fundRows = @(
    @("163801", "中银中国混合(LOF)", "10.14", "89.19", "6.72", "0.6814", 7),
    @("163805", "中银动态策略混合", "5.68", "89.16", "6.62", "0.3760", 7),
    @("009411", "中银科技创新一年定期开放混合", "3.45", "90.71", "6.77", "0.2336", 7),
    @("163809", "中银蓝筹精选灵活配置混合", "3.36", "79.31", "6.67", "0.2241", 5),
    @("501015", "财通多策略升级混合（LOF）", "2.63", "94.52", "8.45", "0.2222", 3),
    @("000612", "华宝生态中国混合", "6.45", "91.21", "3.36", "0.2167", 8),
    @("001118", "华宝事件驱动混合", "6.29", "92.73", "3.23", "0.2032", 8),
    @("011011", "融通产业趋势精选2年封闭运作混合", "3.07", "94.85", "3.44", "0.1056", 10),
    @("005851", "财通新视野灵活配置混合A", "0.87", "85.27", "7.39", "0.0643", 2),
    @("005959", "财通新视野灵活配置混合C", "0.40", "85.27", "7.39", "0.0296", 2),
    @("350001", "天治财富增长混合", "0.98", "69.00", "2.36", "0.0231", 9)
)

# Text columns (fund code / name) plus the numeric-looking text columns
# (scale / position / ratio / market value) must be stored as text, so we
# force a text number format before writing those values.
$q1.Range("B2:B12").NumberFormat = "@"
$q1.Range("D2:G12").NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q1.Cells.Item($r, 1).Value = $i          # A: index (number)
    $q1.Cells.Item($r, 2).Value = $row[0]     # B: fund code (text)
    $q1.Cells.Item($r, 3).Value = $row[1]     # C: fund name (text)
    $q1.Cells.Item($r, 4).Value = $row[2]     # D: fund scale (text)
    $q1.Cells.Item($r, 5).Value = $row[3]     # E: stock position (text)
    $q1.Cells.Item($r, 6).Value = $row[4]     # F: position ratio (text)
    $q1.Cells.Item($r, 7).Value = $row[5]     # G: market value (text)
    $q1.Cells.Item($r, 8).Value = $row[6]     # H: position rank (number)
}

# ---------------------------------------------------------------------------
# Step 2: rewrite $total ("总计") as the refreshed summary table
# ---------------------------------------------------------------------------

# It is still a verbatim copy of the old summary sheet (A1:D6, 5 data
# rows). Extend the styled index column (s="2") one more row, for the new
# 2022-Q1 row that is being inserted at the top.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 11, 2.38),
    @("2021-Q4", 21, 7.27),
    @("2021-Q3", 7, 2.16),
    @("2021-Q2", 27, 6.42),
    @("2021-Q1", 29, 5.51),
    @("2020-Q4", 9, 11.85)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]

    $total.Cells.Item($r, 1).Value = $i        # A: index (number)
    $total.Cells.Item($r, 2).Value = $row[0]   # B: quarter label (text)
    $total.Cells.Item($r, 3).Value = $row[1]   # C: holding count (number)
    $total.Cells.Item($r, 4).Value = $row[2]   # D: holding value (number)
}

Write-Host "2022-Q1 and 总计 sheets updated"
